# The author opened example.xlsx and edited the description text of the
# "Takamatsu Castle" (高松城) row, changing "...あった城である。" to
# "...あったお城である。" (cell D3 on シート1, which is backed by shared
# string index 7). Everything else in the diff (fileVersion/rupBuild,
# absPath, revisionPtr, window geometry, font-metric-driven column widths
# / row heights, cellXfs renumbering, etc.) is incidental re-save noise
# produced by opening the workbook in a different Excel build, not a
# deliberate edit, so only the real content change is applied here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("シート1")

$ws.Range("D3").Value = "高松城（たかまつじょう）は、日本の香川県高松市玉藻町にあったお城である。"

# Reflect the final cursor position recorded in the saved file.
$ws.Range("D9").Select() | Out-Null
